# Add a new "MatrixNodeCollection" worksheet after the existing sheets,
# mirroring the structure of "MatrixEdgeCollection" (same sheet layout:
# a single header cell A1 = "entries").
$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("MatrixEdgeCollection")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Copy the existing collection sheet to the end of the workbook so the new
# sheet inherits the same sheet properties (outline/page setup, margins, etc.)
$source.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "MatrixNodeCollection"

# Ensure the single cell content matches the target schema sheet.
$newSheet.Range("A1").Value = "entries"
